$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Remove the existing hyperlink before shifting rows around - the engine does
# not move hyperlink anchors when rows are inserted, so we recreate them by
# hand once all of the data is in its final place.
$ws.Range("B2").Hyperlinks.Delete()

# Insert two new rows above the existing "Little Thompson Farm" row so the
# new case studies (Catlin Canal Company, Grand Valley) land at rows 2-3.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Inserting rows copies the bold header formatting down onto the new rows;
# put those rows back to normal (non-bold) body formatting.
$ws.Range("A2:F3").Font.Bold = $false

# --- Row 1 (header) is unchanged ---

# The fill order below (row 3 text, row 2 text, both rows' D column, row 6
# text, row 6 D) mirrors the original authoring session's shared-string
# insertion order.

# --- Row 3: Grand Valley Water Users Association ---
$ws.Range("A3").Value2 = "Grand Valley Water Users Association Conserved Consumptive Use Pilot Project"
$ws.Range("B3").Value2 = "http://www.grandvalleywaterusers.com/uploads/8/2/6/0/82606774/03-01-17_ccupp_projectdevelopment_final.pdf"
$ws.Range("C3").Value2 = "Mesa"

# --- Row 2: Catlin Canal Company ---
$ws.Range("A2").Value2 = "Catlin Canal Company Rotational Land Fallowing-Municipal Leasing Pilot Project"
$ws.Range("B2").Value2 = "https://dnrweblink.state.co.us/cwcb/0/edoc/210320/19%2001%2015%202018%20Annual%20Report%20-%20Catlin%20Pilot%20Project%20FINAL.pdf?searchid=3856cf32-c475-4163-840c-5361fa65041f"
$ws.Range("C2").Value2 = "Otero"

$ws.Range("D3").Value2 = "no"
$ws.Range("D2").Value2 = "no"

$ws.Range("E2").Value2 = -103.71323099999999
$ws.Range("F2").Value2 = 38.009126000000002
$ws.Range("E3").Value2 = -108.75243399999999
$ws.Range("F3").Value2 = 39.197431999999999

# --- Row 4: Little Thompson Farm (pre-existing, shifted down from row 2) ---
$ws.Range("A4").Value2 = "Little Thompson Farm"
$ws.Range("B4").Value2 = "https://www.larimer.org/naturalresources/openlands/acquisitions/little-thompson-farm"
$ws.Range("C4").Value2 = "Larimer"
$ws.Range("D4").Value2 = "yes"
$ws.Range("E4").Value2 = -105.10753200000001
$ws.Range("F4").Value2 = 40.286045999999999

# --- Row 5: Maxwell Farm (pre-existing, shifted down from row 3); "no?" -> "no" ---
$ws.Range("A5").Value2 = "Maxwell Farm"
$ws.Range("C5").Value2 = "Larimer"
$ws.Range("D5").Value2 = "no"

# --- Row 6 (new): Yampa Basin Alternative Agricultural Water Transfer Methods Study ---
$ws.Range("A6").Value2 = "Yampa Basin ALternative Agricultural Water Transfer Methods Study"
$ws.Range("B6").Value2 = "https://dnrweblink.state.co.us/cwcb/0/edoc/199193/Yampa%20-%20NC%20Use%20of%20ATM%20to%20Meet%20Non%20%20Consumpt%20Needs_FINALReport%203-28-14_with%20apps.pdf"
$ws.Range("C6").Value2 = "  "
$ws.Range("D6").Value2 = "no"
$ws.Range("E6").Value2 = -107.198581
$ws.Range("F6").Value2 = 40.502360000000003

# --- Hyperlinks, re-added in final cell positions ---
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.larimer.org/naturalresources/openlands/acquisitions/little-thompson-farm") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "http://www.grandvalleywaterusers.com/uploads/8/2/6/0/82606774/03-01-17_ccupp_projectdevelopment_final.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://dnrweblink.state.co.us/cwcb/0/edoc/199193/Yampa%20-%20NC%20Use%20of%20ATM%20to%20Meet%20Non%20%20Consumpt%20Needs_FINALReport%203-28-14_with%20apps.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://dnrweblink.state.co.us/cwcb/0/edoc/210320/19%2001%2015%202018%20Annual%20Report%20-%20Catlin%20Pilot%20Project%20FINAL.pdf?searchid=3856cf32-c475-4163-840c-5361fa65041f") | Out-Null

# --- Column widths (A and B got wider to fit the longer text) ---
$ws.Columns.Item(1).ColumnWidth = 69.21875
$ws.Columns.Item(2).ColumnWidth = 70.77734375
$ws.Columns.Item(4).ColumnWidth = 12.109375

# --- Selection, matching the final saved state ---
$ws.Range("C1:C1048576").Select()
